$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1771736666666667
$ws.Range("H2").Value = 0.531521
$ws.Range("I2").Value = 0.001594420883686164
$ws.Range("J2").Value = 0.001599392044399656
$ws.Range("M2").Value = 0.1771736666666667
$ws.Range("N2").Value = 0.531521
$ws.Range("O2").Value = 0.001594420883686164
$ws.Range("P2").Value = 0.001599392044399656
$ws.Range("Q2").Value = 0.03139050816011112
$ws.Range("R2").Value = 0.282514573441
$ws.Range("S2").Value = [double]"2.5421779543345681e-06"
$ws.Range("T2").Value = [double]"2.5580549116889115e-06"

$ws.Range("G3").Value = 0.1771736666666667
$ws.Range("H3").Value = 0.531521
$ws.Range("I3").Value = 0.001594420883686164
$ws.Range("J3").Value = 0.001599392044399656
$ws.Range("M3").Value = 0.05375533333333333
$ws.Range("N3").Value = 0.161266
$ws.Range("O3").Value = 0.0004837548812342935
$ws.Range("P3").Value = 0.0004852631550440244
$ws.Range("Q3").Value = 0.009524029509555557
$ws.Range("R3").Value = 0.085716265586
$ws.Range("S3").Value = [double]"7.7130888522507759e-07"
$ws.Range("T3").Value = [double]"7.761260296176894e-07"

$ws.Range("G4").Value = 0.1771736666666667
$ws.Range("H4").Value = 0.531521
$ws.Range("I4").Value = 0.001594420883686164
$ws.Range("J4").Value = 0.001599392044399656
$ws.Range("M4").Value = 55.14858233333334
$ws.Range("N4").Value = 165.445747
$ws.Range("O4").Value = 0.4962930046674685
$ws.Range("P4").Value = 0.4978403704304406
$ws.Range("Q4").Value = 9.770876543465224
$ws.Range("R4").Value = 87.93788889118701
$ws.Range("S4").Value = 0.0007912999310691666
$ws.Range("T4").Value = 0.0007962419278474245

$ws.Range("G5").Value = 0.1771736666666667
$ws.Range("H5").Value = 0.531521
$ws.Range("I5").Value = 0.001594420883686164
$ws.Range("J5").Value = 0.001599392044399656
$ws.Range("M5").Value = 1.0361445
$ws.Range("N5").Value = 2.072289
$ws.Range("O5").Value = 0.009324469377408749
$ws.Range("P5").Value = 0.006235694432199138
$ws.Range("Q5").Value = 0.18357752026150004
$ws.Range("R5").Value = 1.1014651215690001
$ws.Range("S5").Value = [double]"1.4867128704632631e-05"
$ws.Range("T5").Value = [double]"9.9733200661665322e-06"

$ws.Range("G6").Value = 0.1771736666666667
$ws.Range("H6").Value = 0.531521
$ws.Range("I6").Value = 0.001594420883686164
$ws.Range("J6").Value = 0.001599392044399656
$ws.Range("M6").Value = 54.70535899999999
$ws.Range("N6").Value = 164.116077
$ws.Range("O6").Value = 0.4923043501902022
$ws.Range("P6").Value = 0.4938392799379165
$ws.Range("Q6").Value = 9.692349040346333
$ws.Range("R6").Value = 87.231141363117
$ws.Range("S6").Value = 0.0007849403370728049
$ws.Range("T6").Value = 0.0007898426155447583

$ws.Range("G7").Value = 0.05375533333333333
$ws.Range("H7").Value = 0.161266
$ws.Range("I7").Value = 0.0004837548812342935
$ws.Range("J7").Value = 0.0004852631550440244
$ws.Range("M7").Value = 0.1771736666666667
$ws.Range("N7").Value = 0.531521
$ws.Range("O7").Value = 0.001594420883686164
$ws.Range("P7").Value = 0.001599392044399656
$ws.Range("Q7").Value = 0.009524029509555557
$ws.Range("R7").Value = 0.085716265586
$ws.Range("S7").Value = [double]"7.7130888522507759e-07"
$ws.Range("T7").Value = [double]"7.761260296176894e-07"

$ws.Range("G8").Value = 0.05375533333333333
$ws.Range("H8").Value = 0.161266
$ws.Range("I8").Value = 0.0004837548812342935
$ws.Range("J8").Value = 0.0004852631550440244
$ws.Range("M8").Value = 0.05375533333333333
$ws.Range("N8").Value = 0.161266
$ws.Range("O8").Value = 0.0004837548812342935
$ws.Range("P8").Value = 0.0004852631550440244
$ws.Range("Q8").Value = 0.0028896358617777774
$ws.Range("R8").Value = 0.026006722755999998
$ws.Range("S8").Value = [double]"2.3401878511800544e-07"
$ws.Range("T8").Value = [double]"2.3548032964328083e-07"

$ws.Range("G9").Value = 0.05375533333333333
$ws.Range("H9").Value = 0.161266
$ws.Range("I9").Value = 0.0004837548812342935
$ws.Range("J9").Value = 0.0004852631550440244
$ws.Range("M9").Value = 55.14858233333334
$ws.Range("N9").Value = 165.445747
$ws.Range("O9").Value = 0.4962930046674685
$ws.Range("P9").Value = 0.4978403704304406
$ws.Range("Q9").Value = 2.964530426189111
$ws.Range("R9").Value = 26.680773835702002
$ws.Range("S9").Value = 0.0002400841635303219
$ws.Range("T9").Value = 0.00024158358886336144

$ws.Range("G10").Value = 0.05375533333333333
$ws.Range("H10").Value = 0.161266
$ws.Range("I10").Value = 0.0004837548812342935
$ws.Range("J10").Value = 0.0004852631550440244
$ws.Range("M10").Value = 1.0361445
$ws.Range("N10").Value = 2.072289
$ws.Range("O10").Value = 0.009324469377408749
$ws.Range("P10").Value = 0.006235694432199138
$ws.Range("Q10").Value = 0.055698292978999996
$ws.Range("R10").Value = 0.334189757874
$ws.Range("S10").Value = [double]"4.510757576241176e-06"
$ws.Range("T10").Value = [double]"3.0259527540594099e-06"

$ws.Range("G11").Value = 0.05375533333333333
$ws.Range("H11").Value = 0.161266
$ws.Range("I11").Value = 0.0004837548812342935
$ws.Range("J11").Value = 0.0004852631550440244
$ws.Range("M11").Value = 54.70535899999999
$ws.Range("N11").Value = 164.116077
$ws.Range("O11").Value = 0.4923043501902022
$ws.Range("P11").Value = 0.4938392799379165
$ws.Range("Q11").Value = 2.9407048081646656
$ws.Range("R11").Value = 26.466343273481996
$ws.Range("S11").Value = 0.0002381546324573873
$ws.Range("T11").Value = 0.00023964200706734251

$ws.Range("G12").Value = 55.14858233333334
$ws.Range("H12").Value = 165.445747
$ws.Range("I12").Value = 0.4962930046674685
$ws.Range("J12").Value = 0.4978403704304406
$ws.Range("M12").Value = 0.1771736666666667
$ws.Range("N12").Value = 0.531521
$ws.Range("O12").Value = 0.001594420883686164
$ws.Range("P12").Value = 0.001599392044399656
$ws.Range("Q12").Value = 9.770876543465224
$ws.Range("R12").Value = 87.93788889118701
$ws.Range("S12").Value = 0.0007912999310691666
$ws.Range("T12").Value = 0.0007962419278474245

$ws.Range("G13").Value = 55.14858233333334
$ws.Range("H13").Value = 165.445747
$ws.Range("I13").Value = 0.4962930046674685
$ws.Range("J13").Value = 0.4978403704304406
$ws.Range("M13").Value = 0.05375533333333333
$ws.Range("N13").Value = 0.161266
$ws.Range("O13").Value = 0.0004837548812342935
$ws.Range("P13").Value = 0.0004852631550440244
$ws.Range("Q13").Value = 2.964530426189111
$ws.Range("R13").Value = 26.680773835702002
$ws.Range("S13").Value = 0.0002400841635303219
$ws.Range("T13").Value = 0.00024158358886336144

$ws.Range("G14").Value = 55.14858233333334
$ws.Range("H14").Value = 165.445747
$ws.Range("I14").Value = 0.4962930046674685
$ws.Range("J14").Value = 0.4978403704304406
$ws.Range("M14").Value = 55.14858233333334
$ws.Range("N14").Value = 165.445747
$ws.Range("O14").Value = 0.4962930046674685
$ws.Range("P14").Value = 0.4978403704304406
$ws.Range("Q14").Value = 3041.3661333764458
$ws.Range("R14").Value = 27372.295200388013
$ws.Range("S14").Value = 0.24630674648186393
$ws.Range("T14").Value = 0.24784503443031833

$ws.Range("G15").Value = 55.14858233333334
$ws.Range("H15").Value = 165.445747
$ws.Range("I15").Value = 0.4962930046674685
$ws.Range("J15").Value = 0.4978403704304406
$ws.Range("M15").Value = 1.0361445
$ws.Range("N15").Value = 2.072289
$ws.Range("O15").Value = 0.009324469377408749
$ws.Range("P15").Value = 0.006235694432199138
$ws.Range("Q15").Value = 57.14190026748051
$ws.Range("R15").Value = 342.85140160488305
$ws.Range("S15").Value = 0.0046276689242439875
$ws.Range("T15").Value = 0.003104380426017055

$ws.Range("G16").Value = 55.14858233333334
$ws.Range("H16").Value = 165.445747
$ws.Range("I16").Value = 0.4962930046674685
$ws.Range("J16").Value = 0.4978403704304406
$ws.Range("M16").Value = 54.70535899999999
$ws.Range("N16").Value = 164.116077
$ws.Range("O16").Value = 0.4923043501902022
$ws.Range("P16").Value = 0.4938392799379165
$ws.Range("Q16").Value = 3016.922994886057
$ws.Range("R16").Value = 27152.30695397452
$ws.Range("S16").Value = 0.24432720516676107
$ws.Range("T16").Value = 0.2458531300573944

$ws.Range("G17").Value = 1.0361445
$ws.Range("H17").Value = 2.072289
$ws.Range("I17").Value = 0.009324469377408749
$ws.Range("J17").Value = 0.006235694432199138
$ws.Range("M17").Value = 0.1771736666666667
$ws.Range("N17").Value = 0.531521
$ws.Range("O17").Value = 0.001594420883686164
$ws.Range("P17").Value = 0.001599392044399656
$ws.Range("Q17").Value = 0.18357752026150004
$ws.Range("R17").Value = 1.1014651215690001
$ws.Range("S17").Value = [double]"1.4867128704632631e-05"
$ws.Range("T17").Value = [double]"9.9733200661665322e-06"

$ws.Range("G18").Value = 1.0361445
$ws.Range("H18").Value = 2.072289
$ws.Range("I18").Value = 0.009324469377408749
$ws.Range("J18").Value = 0.006235694432199138
$ws.Range("M18").Value = 0.05375533333333333
$ws.Range("N18").Value = 0.161266
$ws.Range("O18").Value = 0.0004837548812342935
$ws.Range("P18").Value = 0.0004852631550440244
$ws.Range("Q18").Value = 0.055698292978999996
$ws.Range("R18").Value = 0.334189757874
$ws.Range("S18").Value = [double]"4.510757576241176e-06"
$ws.Range("T18").Value = [double]"3.0259527540594099e-06"

$ws.Range("G19").Value = 1.0361445
$ws.Range("H19").Value = 2.072289
$ws.Range("I19").Value = 0.009324469377408749
$ws.Range("J19").Value = 0.006235694432199138
$ws.Range("M19").Value = 55.14858233333334
$ws.Range("N19").Value = 165.445747
$ws.Range("O19").Value = 0.4962930046674685
$ws.Range("P19").Value = 0.4978403704304406
$ws.Range("Q19").Value = 57.14190026748051
$ws.Range("R19").Value = 342.85140160488305
$ws.Range("S19").Value = 0.0046276689242439875
$ws.Range("T19").Value = 0.003104380426017055

$ws.Range("G20").Value = 1.0361445
$ws.Range("H20").Value = 2.072289
$ws.Range("I20").Value = 0.009324469377408749
$ws.Range("J20").Value = 0.006235694432199138
$ws.Range("M20").Value = 1.0361445
$ws.Range("N20").Value = 2.072289
$ws.Range("O20").Value = 0.009324469377408749
$ws.Range("P20").Value = 0.006235694432199138
$ws.Range("Q20").Value = 1.07359542488025
$ws.Range("R20").Value = 4.294381699521
$ws.Range("S20").Value = [double]"8.69457291702335e-05"
$ws.Range("T20").Value = [double]"3.8883885051759339e-05"

$ws.Range("G21").Value = 1.0361445
$ws.Range("H21").Value = 2.072289
$ws.Range("I21").Value = 0.009324469377408749
$ws.Range("J21").Value = 0.006235694432199138
$ws.Range("M21").Value = 54.70535899999999
$ws.Range("N21").Value = 164.116077
$ws.Range("O21").Value = 0.4923043501902022
$ws.Range("P21").Value = 0.4938392799379165
$ws.Range("Q21").Value = 56.68265684837549
$ws.Range("R21").Value = 340.095941090253
$ws.Range("S21").Value = 0.004590476837713653
$ws.Range("T21").Value = 0.0030794308483100977

$ws.Range("G22").Value = 54.70535899999999
$ws.Range("H22").Value = 164.116077
$ws.Range("I22").Value = 0.4923043501902022
$ws.Range("J22").Value = 0.4938392799379165
$ws.Range("M22").Value = 0.1771736666666667
$ws.Range("N22").Value = 0.531521
$ws.Range("O22").Value = 0.001594420883686164
$ws.Range("P22").Value = 0.001599392044399656
$ws.Range("Q22").Value = 9.692349040346333
$ws.Range("R22").Value = 87.231141363117
$ws.Range("S22").Value = 0.0007849403370728049
$ws.Range("T22").Value = 0.0007898426155447583

$ws.Range("G23").Value = 54.70535899999999
$ws.Range("H23").Value = 164.116077
$ws.Range("I23").Value = 0.4923043501902022
$ws.Range("J23").Value = 0.4938392799379165
$ws.Range("M23").Value = 0.05375533333333333
$ws.Range("N23").Value = 0.161266
$ws.Range("O23").Value = 0.0004837548812342935
$ws.Range("P23").Value = 0.0004852631550440244
$ws.Range("Q23").Value = 2.9407048081646656
$ws.Range("R23").Value = 26.466343273481996
$ws.Range("S23").Value = 0.0002381546324573873
$ws.Range("T23").Value = 0.00023964200706734251

$ws.Range("G24").Value = 54.70535899999999
$ws.Range("H24").Value = 164.116077
$ws.Range("I24").Value = 0.4923043501902022
$ws.Range("J24").Value = 0.4938392799379165
$ws.Range("M24").Value = 55.14858233333334
$ws.Range("N24").Value = 165.445747
$ws.Range("O24").Value = 0.4962930046674685
$ws.Range("P24").Value = 0.4978403704304406
$ws.Range("Q24").Value = 3016.922994886057
$ws.Range("R24").Value = 27152.30695397452
$ws.Range("S24").Value = 0.24432720516676107
$ws.Range("T24").Value = 0.2458531300573944

$ws.Range("G25").Value = 54.70535899999999
$ws.Range("H25").Value = 164.116077
$ws.Range("I25").Value = 0.4923043501902022
$ws.Range("J25").Value = 0.4938392799379165
$ws.Range("M25").Value = 1.0361445
$ws.Range("N25").Value = 2.072289
$ws.Range("O25").Value = 0.009324469377408749
$ws.Range("P25").Value = 0.006235694432199138
$ws.Range("Q25").Value = 56.68265684837549
$ws.Range("R25").Value = 340.095941090253
$ws.Range("S25").Value = 0.004590476837713653
$ws.Range("T25").Value = 0.0030794308483100977

$ws.Range("G26").Value = 54.70535899999999
$ws.Range("H26").Value = 164.116077
$ws.Range("I26").Value = 0.4923043501902022
$ws.Range("J26").Value = 0.4938392799379165
$ws.Range("M26").Value = 54.70535899999999
$ws.Range("N26").Value = 164.116077
$ws.Range("O26").Value = 0.4923043501902022
$ws.Range("P26").Value = 0.4938392799379165
$ws.Range("Q26").Value = 2992.6763033188795
$ws.Range("R26").Value = 26934.086729869927
$ws.Range("S26").Value = 0.24236357321619723
$ws.Range("T26").Value = 0.24387723440959985
